$d = $word.ActiveDocument

# --- Location 1: table cell run "{m" (inside "{m:v.name}") -------------
# Find the *second* occurrence of the literal "{m" (the first occurrence
# belongs to the "{m:for v | self.eClassifiers}" paragraph, which must
# stay untouched).
$r = $d.Content
$r.Find.Execute("{m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0) | Out-Null
$r.Find.Execute("{m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$para1 = $r.Paragraphs(1).Range
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="0033016C" w:rsidR="00B855B8" w:rsidRDefault="00B855B8"><w:r w:rsidR="0033016C"><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r w:rsidRPr="00B90E9D"><w:rPr><w:color w:themeColor="accent6" w:val="F79646"/></w:rPr><w:t>v</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidRPr="00B90E9D"><w:t xml:space="preserve">name}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para1.InsertXML($xml1) | Out-Null

# --- Location 2: paragraph run "{m:" (inside "{m:endfor}") -------------
$r2 = $d.Content
$r2.Find.Execute("{m:endfor}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$para2 = $r2.Paragraphs(1).Range
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="006F5523"><w:r w:rsidR="0033016C"><w:t>{</w:t></w:r><w:r><w:t>m:</w:t></w:r><w:r><w:t xml:space="preserve">endfor}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para2.InsertXML($xml2) | Out-Null

Write-Output "done"
